$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Clear the "Result" data cells (D2:D5) entirely - values and formatting -
# leaving only the D1 header ("Result") in place.
$ws.Range("D2:D5").Clear()

# Move/collapse the selection to B7 (matches the saved sheetView selection).
$ws.Range("B7").Select()

$wb.Save()
